# Update the "Förändrad" (Changed) date column C for rows 2-6
# from 45224 (2023-10-25) to 45233 (2023-11-03), keeping existing
# cell formatting/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [DateTime]::FromOADate(45233)

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
